$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to text format before writing, so that
# numeric-looking strings (e.g. "1.00", "0.513") keep their exact text
# representation instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "94.961.86"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("D3").Value = "3.466.19"
$ws.Range("E3").Value = "  +4.31%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "240.18"
$ws.Range("E5").Value = "  -3.26%  "

$ws.Range("D6").Value = "644.99"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("E7").Value = "  +6.77%  "

$ws.Range("E8").Value = "  -3.10%  "

$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  +3.33%  "

$ws.Range("D11").Value = "3.465.95"
$ws.Range("E11").Value = "  +4.39%  "

$ws.Range("D12").Value = "0.199"
$ws.Range("E12").Value = "  -3.15%  "

$ws.Range("D13").Value = "42.00"
$ws.Range("E13").Value = "  +4.90%  "

$ws.Range("D14").Value = "6.15"
$ws.Range("E14").Value = "  +1.83%  "

$ws.Range("D15").Value = "94.843.40"
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").Value = "4.118.39"
$ws.Range("E16").Value = "  +4.58%  "

$ws.Range("E17").Value = "  +3.03%  "

$ws.Range("D18").Value = "8.57"
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("D19").Value = "3.466.75"
$ws.Range("E19").Value = "  +4.79%  "

$ws.Range("D20").Value = "17.94"
$ws.Range("E20").Value = "  +5.94%  "

$ws.Range("E21").Value = "  +9.78%  "

$ws.Range("D22").Value = "0.513"
$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("D23").Value = "503.66"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("E24").Value = "  -5.18%  "

$ws.Range("D26").Value = "6.67"
$ws.Range("E26").Value = "  +2.51%  "

$ws.Range("D27").Value = "92.03"
$ws.Range("E27").Value = "  -3.71%  "

$ws.Range("D28").Value = "12.20"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("D29").Value = "3.650.12"
$ws.Range("E29").Value = "  +4.31%  "

$ws.Range("D30").Value = "11.75"
$ws.Range("E30").Value = "  +8.05%  "

$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  +12.74%  "

$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("D34").Value = "0.185"
$ws.Range("E34").Value = "  -1.40%  "

$ws.Range("D35").Value = "31.14"
$ws.Range("E35").Value = "  +11.98%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  +5.07%  "

$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").Value = "529.91"
$ws.Range("E40").Value = "  +5.42%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("D43").Value = "0.929"
$ws.Range("E43").Value = "  +12.85%  "

$ws.Range("D44").Value = "24.12"
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("D45").Value = "5.72"
$ws.Range("E45").Value = "  +5.16%  "

$ws.Range("E46").Value = "  +3.25%  "

$ws.Range("D47").Value = "0.0418"
$ws.Range("E47").Value = "  -2.24%  "

$ws.Range("D48").Value = "3.51"
$ws.Range("E48").Value = "  -3.44%  "

$ws.Range("E49").Value = "  +10.35%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "3.22"
$ws.Range("E50").Value = "  +3.71%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "53.38"
$ws.Range("E51").Value = "  +0.68%  "

# Restore the original (default) cell style for the Price column now that
# the text values have been written, so formatting matches the source file.
$ws.Range("D2:D51").Style = "Normal"
